# Documentation mit bildern Meine Notizen
#
# This trims/updates a few speaker-notes text bodies:
#   - Slide 37 notes: remove the "JSF makes it easy..." / "In JSP, you can
#     use..." paragraphs, leaving just the trailing blank paragraph.
#   - Slide 38 notes: remove the "You can use jQuery, Dojo..." paragraph
#     (and the blank line after it), replacing it with a single paragraph
#     that just contains a tab character, then keep the remaining
#     "JSF has built-in capabilities..." paragraph and the trailing blank
#     paragraph.

$p = $ppt.ActivePresentation

$tab = [char]9
$lf  = [char]10

# --- Slide 37 notes: drop the two leading paragraphs entirely ---------
$s37 = $p.Slides.Item(37)
$notes37 = $s37.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes37.Text = ""

# --- Slide 38 notes: leading paragraph becomes a single tab -----------
$s38 = $p.Slides.Item(38)
$notes38 = $s38.NotesPage.Shapes.Item(2).TextFrame.TextRange
$body38 = "JSF has built-in capabilities for checking that form values are in the required format and for converting from strings to various other data types. If values are missing or in an improper format, the form can be automatically redisplayed with error messages and with the previously entered values maintained."
$notes38.Text = "" + $tab + $lf + $body38 + $lf
